$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$url = "https://www.fiercebiotech.com/biotech/cullinan-pens-700m-pact-bcma-bispecific-pair-another-autoimmune-t-cell-engager"
$title = '<a href="https://www.fiercebiotech.com/biotech/cullinan-pens-700m-pact-bcma-bispecific-pair-another-autoimmune-t-cell-engager" hreflang="en">Cullinan pens $700M pact for BCMA bispecific to pair with another autoimmune T-cell engager</a>'

$ws.Range("A20").Value = $url
$ws.Range("B20").Value = "BCMA"
$ws.Range("C20").Value = $title

$ws.Hyperlinks.Add($ws.Range("A20"), $url)
$ws.Range("A20").Style = $ws.Range("A19").Style
